# Add a follow-up note to the "Conference of Swiss Economists Abroad" row
# (row 3), reminding to recheck the location and dates in fall.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note = $ws.Range("E3")
$note.Value = "recheck location and dates in fall!"

# Center the note horizontally and align it to the top of the cell.
$note.HorizontalAlignment = -4108   # xlCenter
$note.VerticalAlignment = -4160     # xlTop

# Leave the new note cell selected, matching the saved selection state.
[void]$note.Select()
